#
# The deck ships two theme parts (ppt/theme/theme1.xml = "Office Theme",
# ppt/theme/theme2.xml = "Integral" / "Red Violet") but only one slide
# master, which is wired to theme2.xml (the "Integral" design, shown as
# the presentation's current Design). The author switched the applied
# design back to the plain "Office Theme" palette via the Design gallery,
# which swaps the 12 theme colours that are actually rendered.
#
# PowerPoint's ThemeColorScheme.Colors(i).RGB exposes exactly those 12
# colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that order)
# for the presentation's active theme, so drive the colour swap through
# that collection on a slide in the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette = the "Office Theme" colour scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink). PowerPoint's ColorFormat.RGB stores colours
# as a VBA-style long (0xBBGGRR), so each literal below is the byte-swap
# of the OOXML srgbClr hex value named in the comment.
$officeColors = @(
    0x000000,   # 1  dk1      srgbClr 000000
    0xFFFFFF,   # 2  lt1      srgbClr FFFFFF
    0x6A5444,   # 3  dk2      srgbClr 44546A
    0xE6E6E7,   # 4  lt2      srgbClr E7E6E6
    0xD59B5B,   # 5  accent1  srgbClr 5B9BD5
    0x317DED,   # 6  accent2  srgbClr ED7D31
    0xA5A5A5,   # 7  accent3  srgbClr A5A5A5
    0x00C0FF,   # 8  accent4  srgbClr FFC000
    0xC47244,   # 9  accent5  srgbClr 4472C4
    0x47AD70,   # 10 accent6  srgbClr 70AD47
    0xC16305,   # 11 hlink    srgbClr 0563C1
    0x724F95    # 12 folHlink srgbClr 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
